# Add a new column N ("bron" / "Bron") populated with "VAP" for every data
# row, mirroring the sticky header styling already used by columns A-C, J
# in rows 2-3 (style index carrying font 18 / fill 33 = orange FFC000).
#
# Rows 2 and 3 are the "sticky" header rows (customFormat, row style = 5);
# writing directly into column N there lets the cell pick up the row's own
# default style automatically (no explicit formatting call needed - exactly
# what a user gets by typing into an already-formatted row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "bron"
$ws.Range("N3").Value = "Bron"

# Data rows 4..66 all get the literal value "VAP" in column N (14).
for ($r = 4; $r -le 66; $r++) {
    $cell = $ws.Cells.Item($r, 14)
    $cell.Value = "VAP"
    $cell.Interior.Color = 49407
}

# Move the active selection in the frozen bottom-right pane from E2 to A3
# (matches the saved view state after the edit).
$ws.Range("A3").Select()
